$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 165, shifting existing rows 165-175 down to 166-176
$ws.Rows.Item(165).Insert()

# Populate the newly inserted row 165 with the new data record
$ws.Cells.Item(165, 1).Value = 5
$ws.Cells.Item(165, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(165, 3).Value = "Maule"
$ws.Cells.Item(165, 4).Value = 44706
$ws.Cells.Item(165, 5).Value = 7
$ws.Cells.Item(165, 6).Value = 100112017
$ws.Cells.Item(165, 7).Value = "Apio"
$ws.Cells.Item(165, 8).Value = "Americana (o)"
$ws.Cells.Item(165, 9).Value = "Primera"
$ws.Cells.Item(165, 10).Value = 600
$ws.Cells.Item(165, 11).Value = 6500
$ws.Cells.Item(165, 12).Value = 6500
$ws.Cells.Item(165, 13).Value = 6500
$ws.Cells.Item(165, 14).Value = "`$/docena de matas"
$ws.Cells.Item(165, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(165, 16).Value = 1083
$ws.Cells.Item(165, 17).Value = 6
$ws.Cells.Item(165, 18).Value = "Hortaliza"

# Match date style/format of column D (numeric date serial with custom date-time format)
$ws.Cells.Item(165, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
